$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new "Save"
# header cell (H1) so it reuses the same cell style as the rest of row 1.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
